# "read setting-data from CSV!" -- add a second sheet ("Sheet2") that is a
# copy of Sheet1 with the M/E mode-row header filled in on both sheets, and
# the "ON AIR" / "BKGD + KEY" status cells filled in on the new sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 3 header labels (applies to Sheet1 as well) -----------------------
$ws1.Range("B3").Value = "Mode1"
$ws1.Range("C3").Value = "Mode2"
$ws1.Range("F3").Value = "M/E1"
$ws1.Range("G3").Value = "M/E2"
$ws1.Range("H3").Value = "M/E3"
$ws1.Range("I3").Value = "M/E4"

# --- Duplicate Sheet1 -> Sheet2 (keeps formatting/styles identical) -------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# --- Sheet2-only content changes -------------------------------------------
$ws2.Range("A4").ClearContents()

$ws2.Range("B5:E5").ClearContents()
$ws2.Range("F5").Value = "ON AIR"
$ws2.Range("G5").Value = "ON AIR"
$ws2.Range("H5").Value = "ON AIR"
$ws2.Range("I5").Value = "ON AIR"
$ws2.Range("J5").ClearContents()

$ws2.Range("B7:D7").ClearContents()
$ws2.Range("E7").Value = "BKGD"
$ws2.Range("F7").Value = "KEY1"
$ws2.Range("G7").Value = "KEY2"
$ws2.Range("H7").Value = "KEY3"
$ws2.Range("I7").Value = "KEY4"
$ws2.Range("J7").ClearContents()

# --- Selections + active sheet/tab -----------------------------------------
$ws1.Range("D11").Select()
$ws2.Range("F5").Select()
$ws2.Activate()
